$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-07-18T16:07:51-04:00"

# --- Elements sheet: collapse the Extension.value[x] slicing ---
$els = $wb.Worksheets.Item("Elements")

# Row 6 (Extension.value[x]) absorbs the more specific
# Extension.value[x]:valueCodeableConcept slice that used to live on row 7:
# the type list collapses to the single fixed type, short text/binding
# information move up, and the (now unused) slicing metadata is cleared.
$els.Range("K6").Value = "CodeableConcept`n"
$els.Range("L6").Value = "Identidad De Genero"
$els.Range("X6").Value = "required"
$els.Range("Y6").Value = ""
$els.Range("Z6").Value = "https://hl7chile.cl/fhir/ig/clcore/ValueSet/VSIdentidaddeGenero"
$els.Range("AB6").Value = ""
$els.Range("AC6").Value = ""
$els.Range("AE6").Value = ""

# Row 7 (the separate valueCodeableConcept slice row) is now redundant; remove it.
$els.Rows.Item(7).Delete()
